# Additional companies sent for questionaire
# Remove the "Parent company" (col B) and "Location County/City" (col E)
# columns from the Known Locomotive List sheet, shifting remaining
# columns left, and update the selection to reflect the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Parent company") entirely.
$ws.Columns.Item(2).Delete()

# After the above shift, the original column E ("Location County/City")
# is now column D - delete it too.
$ws.Columns.Item(4).Delete()

# Reflect the post-edit selection recorded in the sheet view.
[void]$ws.Range("A2:I2").Select()
